$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 16 and 17 need to exist with the same cell style (column A) as the
# rest of the table; clone the formatting from row 15 before filling values.
$ws.Range("A15").Copy($ws.Range("A16:A17"))

# Final target data (rows 8..17): A, B(name), C(from_bus), D(to_bus), E(in_service)
$final = @{
    8  = @(6,  "line7", 14, 11, $true)
    9  = @(7,  "line8", 16, 9,  $true)
    10 = @(8,  "extr1", 5,  12, $true)
    11 = @(9,  "extr2", 5,  9,  $true)
    12 = @(10, "extr3", 10, 11, $true)
    13 = @(11, "extr4", 7,  8,  $false)
    14 = @(12, "extr5", 9,  11, $false)
    15 = @(13, "extr6", 7,  11, $false)
    16 = @(14, "extr7", 5,  7,  $false)
    17 = @(15, "extr8", 8,  5,  $false)
}

# Write bottom-up so we never clobber data we still need to read,
# and never use Insert() (keeps existing per-cell styles / shared-string table order intact).
foreach ($r in ($final.Keys | Sort-Object -Descending)) {
    $vals = $final[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
}
